# Apply marksheet corrections: Corr/total marks
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# "Marking" row: correct mark value per right answer changes from 3 to 5
$ws.Range("B11").Value = 5

# "Total" row: total correct marks changes from 75 to 125
$ws.Range("B12").Value = 125

# Corr/total marks text changes from "73/84" to "125/140"
$ws.Range("E12").Value = "125/140"
